# Update BunkerPrices at 2025-04-14 02:42
# Adds a new data row (row 27) below the existing last row (26), and
# shifts the "short date" number format that was on the previous last
# row's Date cell (Y26) down onto the new last row's Date cell (Y27).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 27

# Row 27 values, column by column (A .. AV), in the same order as the
# header row (Mumbai, Singapore, Rotterdam, ..., Gothenburg), with the
# "Date" column (Y) holding the serial date 2025-04-10.
$values = @{
    "A"  = 567
    "B"  = 478
    "C"  = 442
    "D"  = 549
    "E"  = 494
    "F"  = 532
    "G"  = 475
    "H"  = 569
    "I"  = 512
    "J"  = 442
    "K"  = 571
    "L"  = 483
    "M"  = 447
    "N"  = 505
    "O"  = 555
    "P"  = 483
    "Q"  = 619
    "R"  = 495
    "S"  = 475
    "T"  = 480
    "U"  = 620
    "V"  = 530
    "W"  = 589
    "X"  = 475
    "Y"  = 45757
    "Z"  = 846
    "AA" = 555
    "AB" = 519.5
    "AC" = 512
    "AD" = 538
    "AE" = 500
    "AF" = 502
    "AG" = 750
    "AH" = 459
    "AI" = 740
    "AJ" = 475
    "AK" = 486
    "AL" = 550
    "AM" = 535
    "AN" = 487
    "AO" = 536
    "AP" = 523
    "AQ" = 563
    "AR" = 545
    "AS" = 620
    "AT" = 632
    "AU" = 489
    "AV" = 470
}

foreach ($col in $values.Keys) {
    $ws.Range("$col$newRow").Value = $values[$col]
}

# The previous last row's Date cell used the "date only" format; now that
# it is no longer the last row it switches to the "date + time" format
# that the rest of the column uses, and the new last row's Date cell
# takes over the "date only" format.
$ws.Range("Y26").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("Y27").NumberFormat = "YYYY-MM-DD"
